# Generate Report for Handoff
# Updates the "Priority" and "Latest Handoff Datetime" columns for the
# handed-back files on both the zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 4-7 (a248446b, a51e6033, 1c61f6b7, 3cc011ab files)
# Latest HO Xliff Generate Date (column G) is refreshed to the new
# de-de handoff generation time (shared with the de-de sheet's value).
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = "2016-11-14 07:57:39"
}

# zh-cn sheet: rows 4-7 (a248446b, a51e6033, 1c61f6b7, 3cc011ab files)
# Priority (column E) moves from "low" to "ht"; Latest Handoff Datetime
# (column H) is refreshed to the new handoff generation time.
foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = "2016-11-14 07:57:25"
}

# de-de sheet: rows 4-7, same change plus its own handoff timestamp.
foreach ($row in 4..7) {
    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = "2016-11-14 07:57:39"
}
